# Update documentation of the DNA-component upload form.
#
# Summary of the change:
#  - The old one-line hint in cell A2 ("Clear the example data ...") is
#    removed from the sheet and replaced by a floating text box that
#    holds a longer, 3-paragraph explanation of how to use the sheet.
#  - Because row 2 disappears and the text box pushes the table down by
#    three more blank rows, the whole data table (previously starting at
#    row 4) now starts at row 7, and the table grows from 42 to 45 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# 1) Remove the old inline hint row ("Clear the example data ...").
$ws.Rows("2").Delete()

# 2) Make room: the header row (old row 4, now row 3) needs to end up at
#    row 7, so insert four fresh blank rows above it.
$ws.Rows("3:6").Insert()

# 3) Add the explanatory text box that replaces the old A2 hint, anchored
#    over the now-empty rows 2-6 (same geometry Excel would use for a box
#    spanning columns A-H, rows 2-6).
$shp = $ws.Shapes.AddTextbox(1, 0.75, 17.25, 834.75, 38.25)
$shp.Name = "TextBox 1"

$line1 = "* Clear the example data (except Status and Type columns) before inserting your own data!"
$line2 = "* Markers can be referenced by their ID or by their Name (but by ID is safer). You do NOT need to specify markers for plasmid constructs -- markers will be taken from the associated Vector Backbnone."
$line3 = "* If left empty, the new construct's name will be automatically composed from insert and vector names (only for plasmid entries)."
$text = $line1 + "`n" + $line2 + "`n" + $line3

$shp.TextFrame.Characters().Text = $text
$shp.TextFrame.Characters().Font.Size = 11
$shp.TextFrame.WordWrap = $true
$shp.TextFrame.AutoSize = $false

try { $shp.Fill.ForeColor.RGB = 16777215 } catch {}
try { $shp.Line.Weight = 0.75 } catch {}
try { $shp.Line.ForeColor.RGB = 12632256 } catch {}

# 4) Move the active selection to where the user last left off (E21).
$ws.Range("E21").Select()
